# "9th Stab- Cosmetic Changes"
# Insert two new weekly columns ("Jun_17" and "Jun_15") to the left of the
# existing "Jun_13" / "Jun_10" columns, shifting the old data right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at B:C. This pushes the existing "Jun_13" column
# (old B) to D and the existing "Jun_10" column (old C) to E, carrying their
# values/styles with them automatically.
$ws.Columns("B:C").Insert()

# --- New column B = "Jun_17" -------------------------------------------
$ws.Range("B1").Value = "Jun_17"
$ws.Range("B2:B27").Value = "UN"

# --- New column C = "Jun_15" --------------------------------------------
$ws.Range("C1").Value = "Jun_15"
$ws.Range("C2:C27").Value = "UN"
$ws.Range("C18").Value = "6/13/2018,Raises Target,Outperform,$70.00"

# Row 18 previously had its highlight fill on column B ("Jun_13"); carry
# that same highlight onto the duplicated value that now also appears in
# the new "Jun_15" column (C18), in addition to the shifted original cell
# (D18, already highlighted by the column insert).
$ws.Range("D18").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Restore explicit column widths (8 chars) on C, D, E ----------------
$ws.Columns("C:E").ColumnWidth = 7.1666667
